$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells: "<field>_old" -> "<field>_FV2410",
#        "<field>_new" -> "<field>_FV2504" (columns A1:J1 and L1:U1; K1 "diff" stays as is) ---
$fv2410 = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)
$fv2504 = @(
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $fv2410.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2410[$i]
}
for ($i = 0; $i -lt $fv2504.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2504[$i]
}

# --- 2. Turn the used range into an Excel Table ("Table1") with a header row ---
$tableRange = $ws.Range("A1:U63")
$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Header renaming, table creation and freeze panes applied"
